$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
$hdr.Range.Font.Bold = $hdr.Range.Font.Bold
